$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "A1"
$ws.Range("E2").Value = "A1"
$ws.Range("F2").Value = "M3"
$ws.Range("H2").Value = "DO"
$ws.Range("J2").Value = "M3"
$ws.Range("K2").Value = "DO"
$ws.Range("O2").Value = "M1"
$ws.Range("S2").Value = "DO"
$ws.Range("T2").Value = "M1"
$ws.Range("X2").Value = "A1"
$ws.Range("Y2").Value = "M3"
$ws.Range("Z2").Value = "DO"
$ws.Range("AA2").Value = "M1"
$ws.Range("D3").Value = "DO"
$ws.Range("E3").Value = "M2"
$ws.Range("H3").Value = "A2"
$ws.Range("I3").Value = "M2"
$ws.Range("J3").Value = "DO"
$ws.Range("K3").Value = "M2"
$ws.Range("L3").Value = "A1"
$ws.Range("M3").Value = "A2"
$ws.Range("O3").Value = "M1"
$ws.Range("P3").Value = "DO"
$ws.Range("Q3").Value = "PH"
$ws.Range("R3").Value = "PH"
$ws.Range("S3").Value = "A2"
$ws.Range("T3").Value = "A2"
$ws.Range("X3").Value = "PH"
$ws.Range("Z3").Value = "A2"
$ws.Range("AA3").Value = "A2"
$ws.Range("AB3").Value = "DO"
$ws.Range("AC3").Value = "M2"
$ws.Range("B4").Value = "DO"
$ws.Range("C4").Value = "M1"
$ws.Range("D4").Value = "M1"
$ws.Range("H4").Value = "M3"
$ws.Range("K4").Value = "A1"
$ws.Range("L4").Value = "M1"
$ws.Range("N4").Value = "M3"
$ws.Range("O4").Value = "M1"
$ws.Range("P4").Value = "M3"
$ws.Range("Q4").Value = "PH"
$ws.Range("R4").Value = "PH"
$ws.Range("S4").Value = "DO"
$ws.Range("T4").Value = "M1"
$ws.Range("U4").Value = "M1"
$ws.Range("W4").Value = "M2"
$ws.Range("X4").Value = "PH"
$ws.Range("Z4").Value = "DO"
$ws.Range("AA4").Value = "M2"
$ws.Range("AB4").Value = "A2"
$ws.Range("AC4").Value = "M2"
$ws.Range("B5").Value = "A2"
$ws.Range("C5").Value = "DO"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = "M2"
$ws.Range("F5").Value = "M1"
$ws.Range("G5").Value = "M2"
$ws.Range("H5").Value = "A1"
$ws.Range("I5").Value = "DO"
$ws.Range("K5").Value = "A2"
$ws.Range("L5").Value = "M1"
$ws.Range("N5").Value = "A1"
$ws.Range("O5").Value = "A2"
$ws.Range("P5").Value = "DO"
$ws.Range("Q5").Value = "PH"
$ws.Range("R5").Value = "PH"
$ws.Range("S5").Value = "M2"
$ws.Range("U5").Value = "A2"
$ws.Range("V5").Value = "M2"
$ws.Range("W5").Value = "DO"
$ws.Range("X5").Value = "PH"
$ws.Range("Z5").Value = "M2"
$ws.Range("AA5").Value = "M2"
$ws.Range("AC5").Value = "A2"
$ws.Range("C6").Value = "M1"
$ws.Range("D6").Value = "A1"
$ws.Range("E6").Value = "M1"
$ws.Range("F6").Value = "A1"
$ws.Range("G6").Value = "M3"
$ws.Range("H6").Value = "M1"
$ws.Range("I6").Value = "DO"
$ws.Range("K6").Value = "A2"
$ws.Range("L6").Value = "A2"
$ws.Range("N6").Value = "A1"
$ws.Range("Q6").Value = "M2"
$ws.Range("R6").Value = "A1"
$ws.Range("S6").Value = "M2"
$ws.Range("T6").Value = "A1"
$ws.Range("V6").Value = "DO"
$ws.Range("W6").Value = "DO"
$ws.Range("X6").Value = "M2"
$ws.Range("Y6").Value = "A2"
$ws.Range("Z6").Value = "A2"
$ws.Range("AA6").Value = "M1"
$ws.Range("AB6").Value = "M1"
$ws.Range("D7").Value = "M3"
$ws.Range("F7").Value = "DO"
$ws.Range("H7").Value = "A1"
$ws.Range("K7").Value = "M1"
$ws.Range("L7").Value = "M3"
$ws.Range("N7").Value = "M1"
$ws.Range("O7").Value = "DO"
$ws.Range("R7").Value = "M3"
$ws.Range("S7").Value = "A1"
$ws.Range("T7").Value = "DO"
$ws.Range("V7").Value = "A1"
$ws.Range("Y7").Value = "M1"
$ws.Range("Z7").Value = "M3"
$ws.Range("AA7").Value = "A1"
$ws.Range("B8").Value = "M2"
$ws.Range("C8").Value = "M2"
$ws.Range("D8").Value = "A2"
$ws.Range("F8").Value = "DO"
$ws.Range("I8").Value = "A2"
$ws.Range("J8").Value = "M2"
$ws.Range("M8").Value = "M2"
$ws.Range("O8").Value = "M1"
$ws.Range("P8").Value = "M2"
$ws.Range("R8").Value = "M1"
$ws.Range("U8").Value = "M2"
$ws.Range("W8").Value = "M2"
$ws.Range("Y8").Value = "M1"
$ws.Range("Z8").Value = "A2"
$ws.Range("AB8").Value = "M1"
$ws.Range("C9").Value = "A2"
$ws.Range("D9").Value = "DO"
$ws.Range("E9").Value = "M2"
$ws.Range("F9").Value = "A2"
$ws.Range("G9").Value = "M1"
$ws.Range("J9").Value = "A1"
$ws.Range("K9").Value = "M2"
$ws.Range("L9").Value = "M1"
$ws.Range("M9").Value = "A2"
$ws.Range("N9").Value = "A2"
$ws.Range("R9").Value = "A1"
$ws.Range("V9").Value = "DO"
$ws.Range("W9").Value = "M1"
$ws.Range("Y9").Value = "A2"
$ws.Range("Z9").Value = "M1"
$ws.Range("AA9").Value = "M2"
$ws.Range("B10").Value = "A1"
$ws.Range("C10").Value = "M2"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = "M2"
$ws.Range("F10").Value = "M1"
$ws.Range("G10").Value = "A2"
$ws.Range("I10").Value = "M2"
$ws.Range("J10").Value = "DO"
$ws.Range("L10").Value = "A1"
$ws.Range("M10").Value = "A2"
$ws.Range("N10").Value = "M1"
$ws.Range("O10").Value = "A2"
$ws.Range("Q10").Value = "PH"
$ws.Range("R10").Value = "PH"
$ws.Range("S10").Value = "M2"
$ws.Range("T10").Value = "A2"
$ws.Range("U10").Value = "DO"
$ws.Range("V10").Value = "A2"
$ws.Range("X10").Value = "PH"
$ws.Range("Z10").Value = "M2"
$ws.Range("AA10").Value = "A2"
$ws.Range("AB10").Value = "M2"
